$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-34
# from serial date 45623 (2024-11-27) to 45624 (2024-11-28).
for ($row = 2; $row -le 34; $row++) {
    $ws.Cells.Item($row, 3).Value = 45624
}
